# Adds 2022 data to the SNMV homicide-by-age sheet: the former "Total"
# column (G) becomes the 2022 data column, and a brand-new "Total" column
# (H) is appended with row/column SUM formulas.
# Commit message: "added new data for 2023 in SNMV"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108

# ---------------------------------------------------------------------
# 1. Column G: was the "Total" column (sum of 2017-2021), now becomes the
#    2022 data column. Column H becomes the new "Total" column.
# ---------------------------------------------------------------------
$ws.Range("G1").Value = 2022
$ws.Range("H1").Value = "Total "

# Corrected 2021 values (column F) and new 2022 values (column G)
$F = @{2=3; 3=1; 4=7; 5=8; 6=2; 7=4; 8=4; 9=2; 10=2; 11=2; 12=2; 13=0}
$G = @{2=1; 3=4; 4=4; 5=10; 6=12; 7=5; 8=3; 9=5; 10=2; 11=3; 12=4; 13=0}

foreach ($r in 2..13) {
    $ws.Cells.Item($r, 6).Value = $F[$r]   # column F
    $ws.Cells.Item($r, 7).Value = $G[$r]   # column G
}

# Row 12 (60 años o más): 2020 value corrected from 4 to 5
$ws.Range("E12").Value = 5

# ---------------------------------------------------------------------
# 2. Column H: new row totals = SUM(B:G) for each age-group row.
# ---------------------------------------------------------------------
$ws.Range("H2").Formula = "=SUM(B2:G2)"
$ws.Range("H3:H13").Formula = "=SUM(B3:G3)"
$ws.Range("H14").Formula = "=SUM(H2:H13)"

# ---------------------------------------------------------------------
# 3. Row 14 (Total row): add SUM formulas for every year column.
# ---------------------------------------------------------------------
$ws.Range("B14").Formula = "=SUM(B2:B13)"
$ws.Range("C14").Formula = "=SUM(C2:C13)"
$ws.Range("D14").Formula = "=SUM(D2:D13)"
$ws.Range("E14").Formula = "=SUM(E2:E13)"
$ws.Range("F14").Formula = "=SUM(F2:F13)"
$ws.Range("G14").Formula = "=SUM(G2:G13)"

# ---------------------------------------------------------------------
# 4. Formatting: the data columns F:H (rows 1-14, excluding the
#    untouched F1) are reformatted to the default font, centered both
#    ways. Apply as a single union range / single property write each,
#    so the workbook doesn't accumulate intermediate style snapshots.
# ---------------------------------------------------------------------
$fmtRange = $ws.Range("F2:G11,H2:H13,F12:H14,G1:H1")
$fmtRange.ClearFormats()
$fmtRange.HorizontalAlignment = $xlCenter
$fmtRange.VerticalAlignment = $xlCenter

# ---------------------------------------------------------------------
# 5. Selection moves to B15 (where the next row of data would go).
# ---------------------------------------------------------------------
$ws.Range("B15").Select()
